$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.816.27"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.398.96"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("E14").Value = "  -2.90%  "
$ws.Range("D15").Value = "2.829.05"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "61.754.53"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "2.396.99"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "319.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "559.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "2.515.27"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "0.0₃0923"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -6.02%  "
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "151.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.01%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "146.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0525"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.585"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0224"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.14%  "
